# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 4576
$ws1.Range("F13").Value = 696
$ws1.Range("F14").Value = 184
$ws1.Range("F15").Value = 982
$ws1.Range("F19").Value = 69
$ws1.Range("F20").Value = 117
$ws1.Range("F22").Value = 3517
$ws1.Range("F23").Value = 5872
$ws1.Range("F29").Value = 3358
$ws1.Range("F30").Value = 361
$ws1.Range("F32").Value = 2472
$ws1.Range("F37").Value = 261
$ws1.Range("F38").Value = 352
$ws1.Range("F39").Value = 126
$ws1.Range("F40").Value = 1012
$ws1.Range("F43").Value = 23

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 4576
$ws4.Range("F14").Value = 696
$ws4.Range("F15").Value = 184
$ws4.Range("F16").Value = 982
$ws4.Range("F20").Value = 69
$ws4.Range("F21").Value = 117
$ws4.Range("F23").Value = 3517
$ws4.Range("F24").Value = 5872
$ws4.Range("F30").Value = 3358
$ws4.Range("F31").Value = 361
$ws4.Range("F33").Value = 2472
$ws4.Range("F38").Value = 261
$ws4.Range("F39").Value = 352
$ws4.Range("F40").Value = 126
$ws4.Range("F41").Value = 1012
$ws4.Range("F44").Value = 23
